$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header / "updated as of" timestamp text
$ws.Range("A1").Value = "Datos actualizados a 14 de Abril de 2020 a las 18:52"

# Country-label swaps (sorted table re-ordering caused by updated case counts)
$ws.Cells.Item(23, 1).Value = "India"
$ws.Cells.Item(24, 1).Value = "Irlanda"
$ws.Cells.Item(25, 1).Value = "Corea del Sur"

$ws.Cells.Item(151, 1).Value = "Tanzania"
$ws.Cells.Item(152, 1).Value = "San Martin (Parte Holandesa)"

$ws.Cells.Item(164, 1).Value = "Sudan"
$ws.Cells.Item(165, 1).Value = "San Martin (Parte Francesa)"
$ws.Cells.Item(166, 1).Value = "Mongolia"
$ws.Cells.Item(167, 1).Value = "Siria"

# Updated numeric data (Casos totales, Nuevos casos, Casos activos, Recuperados,
# Casos criticos, Muertes hoy, Muertes)
$ws.Cells.Item(4, 2).Value = 591285
$ws.Cells.Item(4, 3).Value = 4344
$ws.Cells.Item(4, 5).Value = 529218
$ws.Cells.Item(4, 7).Value = 978
$ws.Cells.Item(4, 8).Value = 24618

$ws.Cells.Item(15, 2).Value = 26206
$ws.Cells.Item(15, 3).Value = 526
$ws.Cells.Item(15, 5).Value = 17391

$ws.Cells.Item(17, 2).Value = 24169
$ws.Cells.Item(17, 3).Value = 739
$ws.Cells.Item(17, 4).Value = 3046
$ws.Cells.Item(17, 5).Value = 19745
$ws.Cells.Item(17, 7).Value = 50
$ws.Cells.Item(17, 8).Value = 1378

$ws.Cells.Item(23, 2).Value = 10941
$ws.Cells.Item(23, 3).Value = 488
$ws.Cells.Item(23, 4).Value = 1295
$ws.Cells.Item(23, 5).Value = 9278
$ws.Cells.Item(23, 6).Value = 0
$ws.Cells.Item(23, 7).Value = 10
$ws.Cells.Item(23, 8).Value = 368

$ws.Cells.Item(24, 2).Value = 10647
$ws.Cells.Item(24, 3).Value = 0
$ws.Cells.Item(24, 4).Value = 25
$ws.Cells.Item(24, 5).Value = 10257
$ws.Cells.Item(24, 6).Value = 194
$ws.Cells.Item(24, 7).Value = 0
$ws.Cells.Item(24, 8).Value = 365

$ws.Cells.Item(25, 2).Value = 10564
$ws.Cells.Item(25, 3).Value = 27
$ws.Cells.Item(25, 4).Value = 7534
$ws.Cells.Item(25, 5).Value = 2808
$ws.Cells.Item(25, 6).Value = 55
$ws.Cells.Item(25, 7).Value = 5
$ws.Cells.Item(25, 8).Value = 222

$ws.Cells.Item(55, 5).Value = 1616
$ws.Cells.Item(55, 7).Value = 4
$ws.Cells.Item(55, 8).Value = 102

$ws.Cells.Item(58, 2).Value = 2070
$ws.Cells.Item(58, 3).Value = 87
$ws.Cells.Item(58, 4).Value = 691
$ws.Cells.Item(58, 5).Value = 1053
$ws.Cells.Item(58, 7).Value = 13
$ws.Cells.Item(58, 8).Value = 326

$ws.Cells.Item(59, 5).Value = 1762
$ws.Cells.Item(59, 7).Value = 3
$ws.Cells.Item(59, 8).Value = 38

$ws.Cells.Item(113, 2).Value = 283
$ws.Cells.Item(113, 3).Value = 9
$ws.Cells.Item(113, 5).Value = 274

$ws.Cells.Item(151, 2).Value = 53
$ws.Cells.Item(151, 3).Value = 4
$ws.Cells.Item(151, 4).Value = 7
$ws.Cells.Item(151, 5).Value = 43
$ws.Cells.Item(151, 6).Value = 0
$ws.Cells.Item(151, 8).Value = 3

$ws.Cells.Item(152, 2).Value = 52
$ws.Cells.Item(152, 3).Value = 2
$ws.Cells.Item(152, 4).Value = 5
$ws.Cells.Item(152, 5).Value = 38
$ws.Cells.Item(152, 6).Value = 2
$ws.Cells.Item(152, 8).Value = 9

$ws.Cells.Item(164, 3).Value = 3
$ws.Cells.Item(164, 4).Value = 4
$ws.Cells.Item(164, 5).Value = 23
$ws.Cells.Item(164, 6).Value = 0
$ws.Cells.Item(164, 7).Value = 1
$ws.Cells.Item(164, 8).Value = 5

$ws.Cells.Item(165, 2).Value = 32
$ws.Cells.Item(165, 3).Value = 0
$ws.Cells.Item(165, 4).Value = 11
$ws.Cells.Item(165, 5).Value = 19
$ws.Cells.Item(165, 6).Value = 5
$ws.Cells.Item(165, 8).Value = 2

$ws.Cells.Item(166, 2).Value = 30
$ws.Cells.Item(166, 3).Value = 13
$ws.Cells.Item(166, 5).Value = 25
$ws.Cells.Item(166, 8).Value = 0

$ws.Cells.Item(167, 3).Value = 4
$ws.Cells.Item(167, 5).Value = 22
$ws.Cells.Item(167, 8).Value = 2
